$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Locate the old "_GoBack" bookmark (currently sitting between
#    "...implementat estan" and " explicats a continuació. (") and
#    remove it, merging the two runs around it back into one run.
# ------------------------------------------------------------------
$old = $d.Content
$found = $old.Find.Execute("Els principis que el present treball ha implementat estan explicats a continuació", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Els principis que el present treball ha implementat estan explicats a continuació", 2)

# ------------------------------------------------------------------
# 2) Insert the new sentence after "-patrons." and move the
#    "_GoBack" bookmark to sit right after it (before the trailing
#    space run).
# ------------------------------------------------------------------
$target = $d.Content
$target.Find.Execute("-patrons.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$target.Collapse(0)
$target.InsertAfter(" Els ovals representen problemes, els quadrats els patrons ")
$target.Collapse(0)

# Zero-length ranges aren't accepted reliably by Bookmarks.Add in this
# runtime, so insert a one-character placeholder, bookmark that
# character, then delete the character again leaving an empty
# (point) bookmark exactly where we want it.
$target.InsertAfter("#")
$d.Bookmarks.Add("_GoBack", $target)
$target.Text = ""

$target.Collapse(0)
$target.InsertAfter(" ")
